$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").Value = "'24.682.94"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  -0.86%  "

$style = $ws.Range("D3").Style
$ws.Range("D3").Value = "'1.678.49"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  -1.36%  "

$style = $ws.Range("D4").Style
$ws.Range("D4").Value = "'0.9995"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  -0.30%  "

$style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'313.50"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -0.41%  "

$style = $ws.Range("D6").Style
$ws.Range("D6").Value = "'0.9999"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -0.28%  "

$style = $ws.Range("D7").Style
$ws.Range("D7").Value = "'0.3936"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  -1.77%  "

$style = $ws.Range("D8").Style
$ws.Range("D8").Value = "'0.3958"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  -2.68%  "

$style = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.9995"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  -0.31%  "

$style = $ws.Range("D10").Style
$ws.Range("D10").Value = "'1.412"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  -3.56%  "

$style = $ws.Range("D11").Style
$ws.Range("D11").Value = "'51.04"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -5.06%  "

$style = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.08674"
$ws.Range("D12").Style = $style

$style = $ws.Range("D13").Style
$ws.Range("D13").Value = "'25.41"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  -1.85%  "

$style = $ws.Range("D14").Style
$ws.Range("D14").Value = "'7.355"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -1.70%  "

$style = $ws.Range("D15").Style
$ws.Range("D15").Value = "'0.00001322"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -1.76%  "

$style = $ws.Range("D16").Style
$ws.Range("D16").Value = "'7.734"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -3.78%  "

$style = $ws.Range("D17").Style
$ws.Range("D17").Value = "'1.673.38"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  -2.34%  "

$style = $ws.Range("D18").Style
$ws.Range("D18").Value = "'94.16"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  -2.51%  "

$style = $ws.Range("D19").Style
$ws.Range("D19").Value = "'0.07026"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  -2.30%  "

$style = $ws.Range("D20").Style
$ws.Range("D20").Value = "'21.32"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  +1.78%  "

$style = $ws.Range("D21").Style
$ws.Range("D21").Value = "'7.103"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  -1.96%  "

$style = $ws.Range("D22").Style
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  -0.21%  "

$style = $ws.Range("D23").Style
$ws.Range("D23").Value = "'13.99"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -3.69%  "

$style = $ws.Range("D24").Style
$ws.Range("D24").Value = "'24.690.75"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  -0.83%  "

$style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'2.365"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +1.31%  "

$style = $ws.Range("D26").Style
$ws.Range("D26").Value = "'2.798"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -2.97%  "

$style = $ws.Range("D27").Style
$ws.Range("D27").Value = "'23.19"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +0.43%  "

$style = $ws.Range("D28").Style
$ws.Range("D28").Value = "'5.867"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -11.52%  "

$style = $ws.Range("D29").Style
$ws.Range("D29").Value = "'160.60"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -1.71%  "

$style = $ws.Range("D30").Style
$ws.Range("D30").Value = "'146.80"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  +2.18%  "

$style = $ws.Range("D31").Style
$ws.Range("D31").Value = "'8.279"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  +0.95%  "

$style = $ws.Range("D32").Style
$ws.Range("D32").Value = "'2.507"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  +10.40%  "

$style = $ws.Range("D33").Style
$ws.Range("D33").Value = "'1.852.05"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  -1.82%  "

$style = $ws.Range("D34").Style
$ws.Range("D34").Value = "'0.03101"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  -2.10%  "

$style = $ws.Range("D35").Style
$ws.Range("D35").Value = "'0.08325"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  -4.68%  "

$style = $ws.Range("D36").Style
$ws.Range("D36").Value = "'6.984"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  -5.69%  "

$style = $ws.Range("D37").Style
$ws.Range("D37").Value = "'0.2818"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -1.77%  "

$style = $ws.Range("D38").Style
$ws.Range("D38").Value = "'0.9903"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  -4.00%  "

$style = $ws.Range("D39").Style
$ws.Range("D39").Value = "'0.09547"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  +1.20%  "

$style = $ws.Range("D40").Style
$ws.Range("D40").Value = "'1.524"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  +3.51%  "

$style = $ws.Range("D41").Style
$ws.Range("D41").Value = "'10.37"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  -4.62%  "

$style = $ws.Range("D42").Style
$ws.Range("D42").Value = "'0.7941"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -6.97%  "

$ws.Range("E43").Value = "  -3.18%  "

$style = $ws.Range("D44").Style
$ws.Range("D44").Value = "'16.64"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  -6.05%  "

$style = $ws.Range("D45").Style
$ws.Range("D45").Value = "'0.7148"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  -4.26%  "

$style = $ws.Range("D46").Style
$ws.Range("D46").Value = "'2.576"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  -4.46%  "

$style = $ws.Range("D47").Style
$ws.Range("D47").Value = "'4.169"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -1.32%  "

$style = $ws.Range("D48").Style
$ws.Range("D48").Value = "'0.08670"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +3.52%  "

$style = $ws.Range("D49").Style
$ws.Range("D49").Value = "'0.9993"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -0.29%  "

$ws.Range("E50").Value = "  -4.67%  "

$style = $ws.Range("D51").Style
$ws.Range("D51").Value = "'137.87"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  -2.22%  "
